$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $b, $c, $d, $e, $f) {
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
}

# Update existing rows 4-9 with refined decimal values
Set-Row 4 14537.29 7220.92 10310.57 3085.55 26013.9
Set-Row 5 15727.38 7827.74 11159.05 3006.45 28447.95
Set-Row 6 16943.6095594659 8348.571170663299 11978.0486015057 3301.8655208363 31299.271219973
Set-Row 7 18051.47 8508.49 12530.23 3666.16 34042.57
Set-Row 8 19732.43 9754.07 13984.22 4262.57 36049.41
Set-Row 9 20884.48 10391.6 14711.65 4681.45 38520.26

# Add new row 10 (2021年) - copy style (bold/border/center) from A9
$ws.Cells.Item(10, 1).Value = "2021年"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
Set-Row 10 23167.3 11585.8 16546.4 4855.9 43081.5

# Add new row 11 (2022年) - copy style (bold/border/center) from A9
$ws.Cells.Item(11, 1).Value = "2022年"
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
Set-Row 11 24646.19 11965.26 17450.63 5024.63 46075.42

$excel.CutCopyMode = 0
